# "ran the model for jan 11"
# Updates the Away Team picks for several earlier games, flips two
# "Beat Vegas?" results, fills in the "Beat Vegas?" column for the
# most recent batch of games (Jan 10), corrects the Jan 10 Min/GSW
# rows, and appends the new Jan 11 games.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Date number format used by column A in this sheet (style index 2).
$dateFmt = "yyyy\-mm\-dd"

# --- Away Team (column C) corrections on earlier dates ---
$ws.Cells.Item(2, 3).Value  = "SAS"
$ws.Cells.Item(5, 3).Value  = "CHI"
$ws.Cells.Item(6, 3).Value  = "CLE"
$ws.Cells.Item(9, 3).Value  = "UTA"
$ws.Cells.Item(12, 3).Value = "SAS"
$ws.Cells.Item(16, 3).Value = "CLE"
$ws.Cells.Item(19, 3).Value = "UTA"
$ws.Cells.Item(22, 3).Value = "CHI"
$ws.Cells.Item(27, 3).Value = "SAS"
$ws.Cells.Item(32, 3).Value = "CLE"
$ws.Cells.Item(37, 3).Value = "CHI"
$ws.Cells.Item(38, 3).Value = "UTA"
$ws.Cells.Item(41, 3).Value = "SAS"
$ws.Cells.Item(42, 3).Value = "CLE"

# --- "Beat Vegas?" corrections on Jan 9 games ---
$ws.Cells.Item(35, 7).Value = "Yes"
$ws.Cells.Item(36, 7).Value = "Yes"

# --- Fill in "Beat Vegas?" for the Jan 10 games (rows 41-45) ---
$ws.Cells.Item(41, 7).Value = "Yes"
$ws.Cells.Item(42, 7).Value = "Yes"
$ws.Cells.Item(43, 7).Value = "Yes"
$ws.Cells.Item(44, 7).Value = "No"
$ws.Cells.Item(45, 7).Value = "No"

# --- Row 46: MIN @ CHI replaces the old BOS @ MIA entry ---
$ws.Cells.Item(46, 2).Value = "MIN"
$ws.Cells.Item(46, 3).Value = "CHI"
$ws.Cells.Item(46, 4).Value = 2.5
$ws.Cells.Item(46, 5).Value = 1.5
$ws.Cells.Item(46, 6).Value = 1
$ws.Cells.Item(46, 7).Value = "Yes"

# --- Row 47: GSW @ TOR replaces the old MIN @ SAS entry ---
$ws.Cells.Item(47, 2).Value = "GSW"
$ws.Cells.Item(47, 3).Value = "TOR"
$ws.Cells.Item(47, 4).Value = 3
$ws.Cells.Item(47, 5).Value = -6.3
$ws.Cells.Item(47, 6).Value = 9.3
$ws.Cells.Item(47, 7).Value = "Yes"

# --- Row 48: now the first Jan 11 game (UTA @ MEM) ---
$ws.Cells.Item(48, 1).Value = 44207
$ws.Cells.Item(48, 1).NumberFormat = $dateFmt
$ws.Cells.Item(48, 2).Value = "UTA"
$ws.Cells.Item(48, 3).Value = "MEM"
$ws.Cells.Item(48, 4).Value = 2.5
$ws.Cells.Item(48, 5).Value = 8.5
$ws.Cells.Item(48, 6).Value = -6

# --- New Jan 11 games, rows 49-53 ---
$ws.Cells.Item(49, 1).Value = 44207
$ws.Cells.Item(49, 1).NumberFormat = $dateFmt
$ws.Cells.Item(49, 2).Value = "WAS"
$ws.Cells.Item(49, 3).Value = "PHO"
$ws.Cells.Item(49, 4).Value = 5.5
$ws.Cells.Item(49, 5).Value = -3.3
$ws.Cells.Item(49, 6).Value = 8.8

$ws.Cells.Item(50, 1).Value = 44207
$ws.Cells.Item(50, 1).NumberFormat = $dateFmt
$ws.Cells.Item(50, 2).Value = "ORL"
$ws.Cells.Item(50, 3).Value = "MIL"
$ws.Cells.Item(50, 4).Value = 10
$ws.Cells.Item(50, 5).Value = 22
$ws.Cells.Item(50, 6).Value = -12

$ws.Cells.Item(51, 1).Value = 44207
$ws.Cells.Item(51, 1).NumberFormat = $dateFmt
$ws.Cells.Item(51, 2).Value = "CHO"
$ws.Cells.Item(51, 3).Value = "NYK"
$ws.Cells.Item(51, 4).Value = -5
$ws.Cells.Item(51, 5).Value = -4.5
$ws.Cells.Item(51, 6).Value = -0.5

$ws.Cells.Item(52, 1).Value = 44207
$ws.Cells.Item(52, 1).NumberFormat = $dateFmt
$ws.Cells.Item(52, 2).Value = "ATL"
$ws.Cells.Item(52, 3).Value = "PHI"
$ws.Cells.Item(52, 4).Value = -6
$ws.Cells.Item(52, 5).Value = 15.9
$ws.Cells.Item(52, 6).Value = -21.9

$ws.Cells.Item(53, 1).Value = 44207
$ws.Cells.Item(53, 1).NumberFormat = $dateFmt
$ws.Cells.Item(53, 2).Value = "POR"
$ws.Cells.Item(53, 3).Value = "TOR"
$ws.Cells.Item(53, 4).Value = -5
$ws.Cells.Item(53, 5).Value = -2.8
$ws.Cells.Item(53, 6).Value = -2.2
